# "take minimum in account" - add a minimum-guarantee computation to the
# revenue template: per-representation minima, the resulting guaranteed
# total, the shortfall (and its VAT), and rework the Theatre/Company
# split to use that guaranteed total instead of a flat percentage split.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ---------------------------------------------------------------------
# 1. New named ranges backing the minimum-guarantee computation
# ---------------------------------------------------------------------
$wb.Names.Add('MINIMUM_GARANTI_THEATRE_PAR_REPRESENTATION', '=Feuil1!$G$10')
$wb.Names.Add('MINIMUM_GARANTI_TOTAL', '=Feuil1!$G$11')
$wb.Names.Add('MINIMUM_COMPAGNIE_PAR_REPRESENTATION', '=Feuil1!$G$12')
$wb.Names.Add('DIFFERENCE_RECETTE_MINIMUM_GARANTI', '=Feuil1!$G$13')
$wb.Names.Add('TVA_DIFFERENCE_MINIMUM_GARANTI', '=Feuil1!$G$14')
$wb.Names.Add('TOTAL_THEATRE', '=Feuil1!$G$15')
$wb.Names.Add('NOMBRE_REPRESENTATIONS', '=Feuil1!$A$9:$A$9')

# ---------------------------------------------------------------------
# 2. Updated number of entries for the (only) performance
# ---------------------------------------------------------------------
$ws.Range("C9").Value = 520

# ---------------------------------------------------------------------
# 3. New / relabelled captions in column E
# ---------------------------------------------------------------------
$ws.Range("E10").Value = "minimum garanti par rep :"
$ws.Range("E11").Value = "minimum garanti total :"
$ws.Range("E12").Value = "minimum compagnie par rep :"

# ---------------------------------------------------------------------
# 4. New figures in column G (minimum guarantee block)
# ---------------------------------------------------------------------
$ws.Range("G10").Formula = "=100"
$ws.Range("G11").Formula = '=MINIMUM_GARANTI_THEATRE_PAR_REPRESENTATION*SUM(NOMBRE_REPRESENTATIONS)'
$ws.Range("G12").Value = 100
$ws.Range("G13").Formula = '=-MIN(RECETTE_NETTE - MINIMUM_GARANTI_TOTAL , 0)'
$ws.Range("G14").Formula = '=0.2 * DIFFERENCE_RECETTE_MINIMUM_GARANTI'

# ---------------------------------------------------------------------
# 5. Theatre / Company split now takes the guarantee into account
# ---------------------------------------------------------------------
$ws.Range("F15").Value = 0.5
$ws.Range("F16").Value = 0.5
$ws.Range("G15").Formula = '=MINIMUM_GARANTI_TOTAL + TVA_DIFFERENCE_MINIMUM_GARANTI + (RECETTE_NETTE - MIN(RECETTE_NETTE - MINIMUM_GARANTI_TOTAL, SUM(NOMBRE_REPRESENTATIONS) * MINIMUM_COMPAGNIE_PAR_REPRESENTATION)-MINIMUM_GARANTI_TOTAL) * PART_THEATRE'
$ws.Range("G16").Formula = '=RECETTE_NETTE - TOTAL_THEATRE'

# ---------------------------------------------------------------------
# 6. Cosmetic touch-ups that came along with the edit
# ---------------------------------------------------------------------
$ws.Range("E1").ColumnWidth = 33.72
$ws.Range("E12").Select()
